$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the developer vacation-interval strings (shared strings content) ---
# C2 (Alice, Team Alpha): append a second vacation interval
$ws.Range("C2").Value = "2025-03-10;2025-03-15|2025-05-10;2025-05-15"
# C4 (Charlie, Team Beta): trailing separator added
$ws.Range("C4").Value = "2025-04-05;2025-04-10|"
# C5 (Dave, Team Beta): trailing separator added
$ws.Range("C5").Value = "2025-02-25;2025-02-28|"

# --- Turn the A1:D5 range into a real Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:D5"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"

# --- Resize the columns to fit the new table layout ---
$ws.Columns.Item(1).ColumnWidth = 13.21875
$ws.Columns.Item(2).ColumnWidth = 17.5546875
$ws.Columns.Item(4).ColumnWidth = 17.77734375

# --- Move the active selection to C3 ---
$ws.Range("C3").Select() | Out-Null

Write-Output "Done"
